$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.200.82"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.860.18"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  +0.29%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "239.72"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.37%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "42.32"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +8.08%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D12").Value = "2.128.00"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "1.857.13"
$ws.Range("E14").Value = "  +1.78%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.678"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +1.57%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "4.73"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "35.162.37"
$ws.Range("E17").Value = "  +1.00%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "69.90"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +1.55%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "240.72"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +0.36%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "12.24"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +0.76%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.75"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  +0.42%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.30%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "168.70"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -1.65%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.92"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +27.51%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.01"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +3.44%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "17.68"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +1.36%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.01"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("E33").Value = "  +27.23%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.02"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +10.44%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.818"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +16.99%  "
$ws.Range("E37").Value = "  +7.98%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.10"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +5.24%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0201"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +4.10%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "90.14"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").Value = "1.344.47"
$ws.Range("E41").Value = "  +0.24%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0593"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +14.02%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "14.93"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("E45").Value = "  -0.03%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "12.40"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +44.34%  "
$ws.Range("E47").Value = "  -0.49%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "6.61"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +5.48%  "
$ws.Range("D49").Value = "2.044.60"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("E51").Value = "  +0.40%  "
